$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.04684866123356
$ws.Range("C2").Value = 14.85045676676299
$ws.Range("D2").Value = 7.345983262333876
$ws.Range("F2").Value = 43.20575197441379
$ws.Range("G2").Value = 3.697102544084863
$ws.Range("I2").Value = 27.65027656396731
$ws.Range("J2").Value = 10.56280034462994
$ws.Range("M2").Value = 20.16515310504352
$ws.Range("N2").Value = 19.51121820688275
$ws.Range("B3").Value = 17.60516066014521
$ws.Range("C3").Value = 14.41943652876106
$ws.Range("D3").Value = 7.348090612391682
$ws.Range("F3").Value = 43.06383683378426
$ws.Range("G3").Value = 3.700997358722977
$ws.Range("I3").Value = 27.66606152737132
$ws.Range("J3").Value = 10.58654027939648
$ws.Range("M3").Value = 20.04308826280217
$ws.Range("N3").Value = 19.58319746836055
$ws.Range("B4").Value = 17.33361034092493
$ws.Range("C4").Value = 14.15298478268926
$ws.Range("D4").Value = 7.349816688739533
$ws.Range("F4").Value = 42.98892397435331
$ws.Range("G4").Value = 3.703512157776867
$ws.Range("I4").Value = 27.68275153842992
$ws.Range("J4").Value = 10.60238459007639
$ws.Range("M4").Value = 19.97258538269163
$ws.Range("N4").Value = 19.62937420225289
$ws.Range("B5").Value = 17.22304223857712
$ws.Range("C5").Value = 14.04413673013127
$ws.Range("D5").Value = 7.350629178820574
$ws.Range("F5").Value = 42.96148297509298
$ws.Range("G5").Value = 3.704568100674881
$ws.Range("I5").Value = 27.69130589235916
$ws.Range("J5").Value = 10.60916022875889
$ws.Range("M5").Value = 19.94499572360947
$ws.Range("N5").Value = 19.6486910981575
$ws.Range("B6").Value = 17.20469332994862
$ws.Range("C6").Value = 14.02605212143474
$ws.Range("D6").Value = 7.350770694115977
$ws.Range("F6").Value = 42.95711313835706
$ws.Range("G6").Value = 3.70474532360838
$ws.Range("I6").Value = 27.69283200376372
$ws.Range("J6").Value = 10.61030458626776
$ws.Range("M6").Value = 19.94048402118224
$ws.Range("N6").Value = 19.65192886155184
$ws.Range("B7").Value = 17.33211857179671
$ws.Range("C7").Value = 14.15151764630714
$ws.Range("D7").Value = 7.349827203938617
$ws.Range("F7").Value = 42.98854138414158
$ws.Range("G7").Value = 3.703526272342302
$ws.Range("I7").Value = 27.68285981707391
$ws.Range("J7").Value = 10.60247467711348
$ws.Range("M7").Value = 19.97220865074582
$ws.Range("N7").Value = 19.62963269252261
$ws.Range("B8").Value = 17.8947334655343
$ws.Range("C8").Value = 14.70232429591602
$ws.Range("D8").Value = 7.346620441415899
$ws.Range("F8").Value = 43.15429154221732
$ws.Range("G8").Value = 3.698419943838237
$ws.Range("I8").Value = 27.65426303694614
$ws.Range("J8").Value = 10.57072280139529
$ws.Range("M8").Value = 20.12215684784598
$ws.Range("N8").Value = 19.53562626096285
$ws.Range("B9").Value = 18.98768026865707
$ws.Range("C9").Value = 15.76046526480967
$ws.Range("D9").Value = 7.343742394517965
$ws.Range("F9").Value = 43.57557040533569
$ws.Range("G9").Value = 3.689379827989034
$ws.Range("I9").Value = 27.65397230142385
$ws.Range("J9").Value = 10.51851206526602
$ws.Range("M9").Value = 20.45038359277467
$ws.Range("N9").Value = 19.36693800732534
$ws.Range("B10").Value = 19.77483454222781
$ws.Range("C10").Value = 16.51500827006473
$ws.Range("D10").Value = 7.34368261423059
$ws.Range("F10").Value = 43.94250301729141
$ws.Range("G10").Value = 3.683323775956905
$ws.Range("I10").Value = 27.68807233206202
$ws.Range("J10").Value = 10.4862751459793
$ws.Range("M10").Value = 20.71079344996259
$ws.Range("N10").Value = 19.25245994147812
$ws.Range("B11").Value = 20.1276645677153
$ws.Range("C11").Value = 16.85154683329043
$ws.Range("D11").Value = 7.344096371725084
$ws.Range("F11").Value = 44.12155163475432
$ws.Range("G11").Value = 3.680694263642445
$ws.Range("I11").Value = 27.71108138368355
$ws.Range("J11").Value = 10.47293767161396
$ws.Range("M11").Value = 20.83307829697405
$ws.Range("N11").Value = 19.2024156831249
$ws.Range("B12").Value = 20.26037876210178
$ws.Range("C12").Value = 16.97789156397716
$ws.Range("D12").Value = 7.344316003729155
$ws.Range("F12").Value = 44.19106184962479
$ws.Range("G12").Value = 3.679716444136993
$ws.Range("I12").Value = 27.72087424200598
$ws.Range("J12").Value = 10.46807787240232
$ws.Range("M12").Value = 20.87990165937387
$ws.Range("N12").Value = 19.18375613590711
$ws.Range("B13").Value = 20.23183829884211
$ws.Range("C13").Value = 16.95073162038572
$ws.Range("D13").Value = 7.344265909615714
$ws.Range("F13").Value = 44.17601620528161
$ws.Range("G13").Value = 3.679926239903497
$ws.Range("I13").Value = 27.7187171314555
$ws.Range("J13").Value = 10.46911603138601
$ws.Range("M13").Value = 20.8697949634142
$ws.Range("N13").Value = 19.18776187301284
$ws.Range("B14").Value = 20.13860174213265
$ws.Range("C14").Value = 16.86196394257539
$ws.Range("D14").Value = 7.344113181983758
$ws.Range("F14").Value = 44.12723626015403
$ws.Range("G14").Value = 3.680613459362315
$ws.Range("I14").Value = 27.71186540004997
$ws.Range("J14").Value = 10.47253402909527
$ws.Range("M14").Value = 20.8369203191825
$ws.Range("N14").Value = 19.2008747227109
$ws.Range("B15").Value = 20.08137115105964
$ws.Range("C15").Value = 16.80744489348025
$ws.Range("D15").Value = 7.344027816337858
$ws.Range("F15").Value = 44.09757849877683
$ws.Range("G15").Value = 3.681036731397863
$ws.Range("I15").Value = 27.707809178003
$ws.Range("J15").Value = 10.47465249949845
$ws.Range("M15").Value = 20.81684995028019
$ws.Range("N15").Value = 19.20894460201889
$ws.Range("B16").Value = 19.75165841708903
$ws.Range("C16").Value = 16.49286826892359
$ws.Range("D16").Value = 7.343664408289159
$ws.Range("F16").Value = 43.93104285885539
$ws.Range("G16").Value = 3.68349813462857
$ws.Range("I16").Value = 27.68671966874314
$ws.Range("J16").Value = 10.48717348527778
$ws.Range("M16").Value = 20.70287597784874
$ws.Range("N16").Value = 19.25577125143469
$ws.Range("B17").Value = 19.54794266158157
$ws.Range("C17").Value = 16.29807096765841
$ws.Range("D17").Value = 7.343554083533483
$ws.Range("F17").Value = 43.83196084469519
$ws.Range("G17").Value = 3.685040165946479
$ws.Range("I17").Value = 27.67570365687989
$ws.Range("J17").Value = 10.49519461540513
$ws.Range("M17").Value = 20.63391264498476
$ws.Range("N17").Value = 19.2850176185898
$ws.Range("B18").Value = 19.43028424348766
$ws.Range("C18").Value = 16.18540444423022
$ws.Range("D18").Value = 7.343532159751278
$ws.Range("F18").Value = 43.77611550740738
$ws.Range("G18").Value = 3.685938913184761
$ws.Range("I18").Value = 27.67007312900538
$ws.Range("J18").Value = 10.49993308822525
$ws.Range("M18").Value = 20.59460865044684
$ws.Range("N18").Value = 19.30203071963318
$ws.Range("B19").Value = 19.39036800475938
$ws.Range("C19").Value = 16.14715446319463
$ws.Range("D19").Value = 7.343531883780996
$ws.Range("F19").Value = 43.75740475512049
$ws.Range("G19").Value = 3.686245245719326
$ws.Range("I19").Value = 27.66828784271399
$ws.Range("J19").Value = 10.50155891209593
$ws.Range("M19").Value = 20.58136413832349
$ws.Range("N19").Value = 19.30782396494943
$ws.Range("B20").Value = 19.56967995010076
$ws.Range("C20").Value = 16.31887306849269
$ws.Range("D20").Value = 7.34356153267739
$ws.Range("F20").Value = 43.84239014115928
$ws.Range("G20").Value = 3.684874792457936
$ws.Range("I20").Value = 27.67680328781725
$ws.Range("J20").Value = 10.49432782287053
$ws.Range("M20").Value = 20.64121669135171
$ws.Range("N20").Value = 19.28188449565321
$ws.Range("B21").Value = 20.16601292723828
$ws.Range("C21").Value = 16.88806786250046
$ws.Range("D21").Value = 7.344156336985027
$ws.Range("F21").Value = 44.14151804901982
$ws.Range("G21").Value = 3.68041112092626
$ws.Range("I21").Value = 27.71384860421147
$ws.Range("J21").Value = 10.47152490207111
$ws.Range("M21").Value = 20.8465626423981
$ws.Range("N21").Value = 19.19701526866293
$ws.Range("B22").Value = 20.55047765196444
$ws.Range("C22").Value = 17.25363442481862
$ws.Range("D22").Value = 7.344911787592815
$ws.Range("F22").Value = 44.34695430857133
$ws.Range("G22").Value = 3.677598250437745
$ws.Range("I22").Value = 27.74435402793238
$ws.Range("J22").Value = 10.45773405751107
$ws.Range("M22").Value = 20.98376415900585
$ws.Range("N22").Value = 19.1432447268354
$ws.Range("B23").Value = 20.3458075430743
$ws.Range("C23").Value = 17.05915347236055
$ws.Range("D23").Value = 7.344475187422147
$ws.Range("F23").Value = 44.2364123777302
$ws.Range("G23").Value = 3.67909001790654
$ws.Range("I23").Value = 27.72749646836585
$ws.Range("J23").Value = 10.46499274106053
$ws.Range("M23").Value = 20.91027412003622
$ws.Range("N23").Value = 19.17178822303286
$ws.Range("B24").Value = 19.55985419170499
$ws.Range("C24").Value = 16.30947053100475
$ws.Range("D24").Value = 7.343558035611109
$ws.Range("F24").Value = 43.83767157051124
$ws.Range("G24").Value = 3.684949519763277
$ws.Range("I24").Value = 27.67630395599055
$ws.Range("J24").Value = 10.49471930413851
$ws.Range("M24").Value = 20.63791345818767
$ws.Range("N24").Value = 19.28330036068479
$ws.Range("B25").Value = 18.69413542188243
$ws.Range("C25").Value = 15.47758894616669
$ws.Range("D25").Value = 7.344158374057313
$ws.Range("F25").Value = 43.45141380965337
$ws.Range("G25").Value = 3.69172200710811
$ws.Range("I25").Value = 27.64804329510254
$ws.Range("J25").Value = 10.53156079121939
$ws.Range("M25").Value = 20.35808709529748
$ws.Range("N25").Value = 19.41090555064747
